$wb = $excel.ActiveWorkbook

# --- ALC sheet ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H42").Value = 2956.875
$ws.Range("I42").Value = 23.444445
$ws.Range("K42").Value = 70.33333500000001
$ws.Range("M42").Value = 159.666665
$ws.Range("H47").Value = 12901.5
$ws.Range("J47").Value = 17035.334
$ws.Range("L47").Value = 17035.334
$ws.Range("N47").Value = -18979.334
$ws.Range("H69").Value = 8534
$ws.Range("J69").Value = 12489.5
$ws.Range("L69").Value = 37468.5
$ws.Range("N69").Value = -39216.5
$ws.Range("H72").Value = 8534
$ws.Range("J72").Value = 12489.5
$ws.Range("L72").Value = 112405.5
$ws.Range("N72").Value = -121141.5
$ws.Range("H88").Value = 776181.3
$ws.Range("I88").Value = 6536.3335
$ws.Range("K88").Value = 6536.3335
$ws.Range("M88").Value = -6130.3335
$ws.Range("H91").Value = 776181.3
$ws.Range("I91").Value = 6536.3335
$ws.Range("K91").Value = 6536.3335
$ws.Range("M91").Value = -5132.3335
$ws.Range("H99").Value = 5403.375
$ws.Range("J99").Value = 9187
$ws.Range("L99").Value = 27561
$ws.Range("N99").Value = -30557
$ws.Range("H141").Value = 2836.5483
$ws.Range("I141").Value = 2925.4138
$ws.Range("K141").Value = 8776.241399999999
$ws.Range("M141").Value = -3596.241399999999

# --- ARM sheet ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 10794.2
$ws.Range("I28").Value = 10992.75
$ws.Range("J28").Value = 10000
$ws.Range("K28").Value = 10992.75
$ws.Range("L28").Value = 10000
$ws.Range("M28").Value = -10800.75
$ws.Range("N28").Value = -10384
$ws.Range("H45").Value = 1580.3871
$ws.Range("I45").Value = 1392.909
$ws.Range("J45").Value = 2038.6666
$ws.Range("K45").Value = 1392.909
$ws.Range("L45").Value = 2038.6666
$ws.Range("M45").Value = -1015.909
$ws.Range("N45").Value = -2792.6666
$ws.Range("H74").Value = 3539.2593
$ws.Range("I74").Value = 2805.1333
$ws.Range("K74").Value = 2805.1333
$ws.Range("M74").Value = -1931.1333
$ws.Range("H77").Value = 3539.2593
$ws.Range("I77").Value = 2805.1333
$ws.Range("K77").Value = 14025.6665
$ws.Range("M77").Value = -9657.666499999999
$ws.Range("H88").Value = 2037.7858
$ws.Range("I88").Value = 1967.5
$ws.Range("J88").Value = 2065.9
$ws.Range("K88").Value = 1967.5
$ws.Range("L88").Value = 2065.9
$ws.Range("M88").Value = -1561.5
$ws.Range("N88").Value = -2877.9
$ws.Range("H91").Value = 2037.7858
$ws.Range("I91").Value = 1967.5
$ws.Range("J91").Value = 2065.9
$ws.Range("K91").Value = 1967.5
$ws.Range("L91").Value = 2065.9
$ws.Range("M91").Value = -563.5
$ws.Range("N91").Value = -4873.9
$ws.Range("H97").Value = 472.42856
$ws.Range("I97").Value = 319.58334
$ws.Range("J97").Value = 676.2222
$ws.Range("K97").Value = 319.58334
$ws.Range("L97").Value = 676.2222
$ws.Range("M97").Value = 176.41666
$ws.Range("N97").Value = -1668.2222
$ws.Range("H99").Value = 10794.2
$ws.Range("I99").Value = 10992.75
$ws.Range("J99").Value = 10000
$ws.Range("K99").Value = 10992.75
$ws.Range("L99").Value = 10000
$ws.Range("M99").Value = -7997.75
$ws.Range("N99").Value = -15990
$ws.Range("H122").Value = 3338.9744
$ws.Range("I122").Value = 3178.138
$ws.Range("J122").Value = 3805.4
$ws.Range("K122").Value = 9534.414000000001
$ws.Range("L122").Value = 11416.2
$ws.Range("M122").Value = -7084.414000000001
$ws.Range("N122").Value = -16316.2

# --- BSM sheet ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1731.3939
$ws.Range("J86").Value = 2010
$ws.Range("L86").Value = 2010
$ws.Range("N86").Value = -4256
$ws.Range("H89").Value = 1731.3939
$ws.Range("J89").Value = 2010
$ws.Range("L89").Value = 10050
$ws.Range("N89").Value = -21282

# --- CRP sheet ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 3594.4055
$ws.Range("I99").Value = 3061
$ws.Range("K99").Value = 3061
$ws.Range("M99").Value = -1563
$ws.Range("H126").Value = 3594.4055
$ws.Range("I126").Value = 3061
$ws.Range("K126").Value = 9183
$ws.Range("M126").Value = -6713

# --- CUL sheet ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H59").Value = 6827.1816
$ws.Range("I59").Value = 3749.5
$ws.Range("J59").Value = 7511.1113
$ws.Range("K59").Value = 11248.5
$ws.Range("L59").Value = 22533.3339
$ws.Range("M59").Value = -10708.5
$ws.Range("N59").Value = -23613.3339
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()
$ws.Range("H69").Value = 6457.143
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()
$ws.Range("H72").Value = 6457.143
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()
$ws.Range("H93").Value = 750
$ws.Range("I93").Value = 750
$ws.Range("K93").Value = 2250
$ws.Range("M93").Value = -378

# --- GSM sheet ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H31").Value = 4010.5
$ws.Range("I31").Value = 4010.5
$ws.Range("K31").Value = 4010.5
$ws.Range("M31").Value = -3718.5
$ws.Range("H37").Value = 4010.5
$ws.Range("I37").Value = 4010.5
$ws.Range("K37").Value = 4010.5
$ws.Range("M37").Value = -3733.5
$ws.Range("H70").Value = 24076.936
$ws.Range("J70").Value = 9155.25
$ws.Range("L70").Value = 9155.25
$ws.Range("N70").Value = -9695.25
$ws.Range("H73").Value = 24076.936
$ws.Range("J73").Value = 9155.25
$ws.Range("L73").Value = 9155.25
$ws.Range("N73").Value = -11027.25
$ws.Range("H102").Value = 2633.7646
$ws.Range("I102").Value = 2383.6155
$ws.Range("J102").Value = 3446.75
$ws.Range("K102").Value = 2383.6155
$ws.Range("L102").Value = 3446.75
$ws.Range("M102").Value = -761.6154999999999
$ws.Range("N102").Value = -6690.75
$ws.Range("H122").Value = 3406.2778
$ws.Range("I122").Value = 2093.3572
$ws.Range("J122").Value = 8001.5
$ws.Range("K122").Value = 6280.071599999999
$ws.Range("L122").Value = 24004.5
$ws.Range("M122").Value = -3830.071599999999
$ws.Range("N122").Value = -28904.5
$ws.Range("H126").Value = 5101.1
$ws.Range("I126").Value = 4799.8
$ws.Range("J126").Value = 5402.4
$ws.Range("K126").Value = 14399.4
$ws.Range("L126").Value = 16207.2
$ws.Range("M126").Value = -11929.4
$ws.Range("N126").Value = -21147.2

# --- LTW sheet ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 3903.3333
$ws.Range("I22").Value = 1748.6
$ws.Range("K22").Value = 1748.6
$ws.Range("M22").Value = -1453.6
$ws.Range("H27").Value = 3903.3333
$ws.Range("I27").Value = 1748.6
$ws.Range("K27").Value = 1748.6
$ws.Range("M27").Value = -1641.6
$ws.Range("H68").Value = 1415
$ws.Range("I68").Value = 1869.3077
$ws.Range("J68").Value = 1086.8889
$ws.Range("K68").Value = 1869.3077
$ws.Range("L68").Value = 1086.8889
$ws.Range("M68").Value = -1120.3077
$ws.Range("N68").Value = -2584.8889
$ws.Range("H71").Value = 1415
$ws.Range("I71").Value = 1869.3077
$ws.Range("J71").Value = 1086.8889
$ws.Range("K71").Value = 9346.538500000001
$ws.Range("L71").Value = 5434.4445
$ws.Range("M71").Value = -5602.538500000001
$ws.Range("N71").Value = -12922.4445

# --- WVR sheet ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 2390.611
$ws.Range("I81").Value = 1083.4615
$ws.Range("J81").Value = 5789.2
$ws.Range("K81").Value = 2166.923
$ws.Range("L81").Value = 11578.4
$ws.Range("M81").Value = -1105.923
$ws.Range("N81").Value = -13700.4
$ws.Range("H84").Value = 2390.611
$ws.Range("I84").Value = 1083.4615
$ws.Range("J84").Value = 5789.2
$ws.Range("K84").Value = 10834.615
$ws.Range("L84").Value = 57892
$ws.Range("M84").Value = -5530.614999999998
$ws.Range("N84").Value = -68500
$ws.Range("H122").Value = 6352.579
$ws.Range("I122").Value = 6317.846
$ws.Range("J122").Value = 6427.8335
$ws.Range("K122").Value = 18953.538
$ws.Range("L122").Value = 19283.5005
$ws.Range("M122").Value = -16503.538
$ws.Range("N122").Value = -24183.5005
$ws.Range("H132").Value = 705.5789
$ws.Range("I132").Value = 612
$ws.Range("J132").Value = 2390
$ws.Range("K132").Value = 1836
$ws.Range("L132").Value = 7170
$ws.Range("M132").Value = 694
$ws.Range("N132").Value = -12230
